$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.582.72"
$ws.Range("E2").Value = "  +9.93%  "
$ws.Range("D3").Value = "2.592.96"
$ws.Range("E3").Value = "  +10.82%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "513.45"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +8.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.05"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +9.99%  "
$ws.Range("E7").Value = "  -0.79%  "
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("D9").Value = "2.646.82"
$ws.Range("E9").Value = "  +13.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.15"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +13.38%  "
$ws.Range("E11").Value = "  +10.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.348"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +6.89%  "
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("D14").Value = "3.042.46"
$ws.Range("E14").Value = "  +10.97%  "
$ws.Range("D15").Value = "60.277.40"
$ws.Range("E15").Value = "  +9.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.29"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +11.80%  "
$ws.Range("E17").Value = "  +9.51%  "
$ws.Range("D18").Value = "2.626.36"
$ws.Range("E18").Value = "  +12.49%  "
$ws.Range("E19").Value = "  +6.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "344.45"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +9.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.59"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +11.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.18"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +10.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.996"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.70"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.425"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +8.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.170"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +12.59%  "
$ws.Range("D27").Value = "2.697.02"
$ws.Range("E27").Value = "  +10.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.989"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.16%  "
$ws.Range("D29").Value = "0.0₃0865"
$ws.Range("E29").Value = "  +17.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.56"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +7.46%  "
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.71"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +9.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "157.25"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +8.10%  "
$ws.Range("E34").Value = "  +8.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.60"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +10.49%  "
$ws.Range("E36").Value = "  +11.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.99"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +10.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.876"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +8.44%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.49"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +13.17%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "305.06"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +22.65%  "
$ws.Range("E41").Value = "  +11.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "35.64"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.16%  "
$ws.Range("E43").Value = "  +11.57%  "
$ws.Range("B44").Value = "SuiNetwork"
$ws.Range("C44").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.804"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +30.53%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0578"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +12.47%  "
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.02"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +20.93%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.02"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +15.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.988"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0239"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +8.55%  "
$ws.Range("D51").Value = "2.009.96"
$ws.Range("E51").Value = "  +12.60%  "
